$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 17:20"

# Update country rows: name (col A) and stats (cols B-H)
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 124356
$ws.Range("C4").Value = 778
$ws.Range("D4").Value = 3238
$ws.Range("E4").Value = 118882
$ws.Range("F4").Value = 2666
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 2236

$ws.Range("A12").Value = "Suiza"
$ws.Range("B12").Value = 14829
$ws.Range("C12").Value = 753
$ws.Range("D12").Value = 1595
$ws.Range("E12").Value = 12934
$ws.Range("F12").Value = 301
$ws.Range("G12").Value = 36
$ws.Range("H12").Value = 300

$ws.Range("A16").Value = "Austria"
$ws.Range("B16").Value = 8648
$ws.Range("C16").Value = 377
$ws.Range("D16").Value = 479
$ws.Range("E16").Value = 8083
$ws.Range("F16").Value = 187
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 86

$ws.Range("A19").Value = "Canada"
$ws.Range("B19").Value = 5655
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 508
$ws.Range("E19").Value = 5084
$ws.Range("F19").Value = 120
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = 63

$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 2139
$ws.Range("C29").Value = 230
$ws.Range("D29").Value = 61
$ws.Range("E29").Value = 2072
$ws.Range("F29").Value = 7
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 6

$ws.Range("A41").Value = "Finlandia"
$ws.Range("B41").Value = 1240
$ws.Range("C41").Value = 73
$ws.Range("D41").Value = 10
$ws.Range("E41").Value = 1219
$ws.Range("F41").Value = 32
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 11

$ws.Range("A44").Value = "India"
$ws.Range("B44").Value = 1024
$ws.Range("C44").Value = 37
$ws.Range("D44").Value = 87
$ws.Range("E44").Value = 910
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 27

$ws.Range("A45").Value = "Islandia"
$ws.Range("B45").Value = 1020
$ws.Range("C45").Value = 57
$ws.Range("D45").Value = 114
$ws.Range("E45").Value = 904
$ws.Range("F45").Value = 19
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 2

$ws.Range("A87").Value = "Republica de Chipre"
$ws.Range("B87").Value = 214
$ws.Range("C87").Value = 35
$ws.Range("D87").Value = 15
$ws.Range("E87").Value = 194
$ws.Range("F87").Value = 3
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 5

$ws.Range("A88").Value = "Albania"
$ws.Range("B88").Value = 212
$ws.Range("C88").Value = 15
$ws.Range("D88").Value = 33
$ws.Range("E88").Value = 169
$ws.Range("F88").Value = 3
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 10

$ws.Range("A89").Value = "Azerbaiyan"
$ws.Range("B89").Value = 209
$ws.Range("C89").Value = 27
$ws.Range("D89").Value = 15
$ws.Range("E89").Value = 190
$ws.Range("F89").Value = 23
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 4

$ws.Range("A90").Value = "Burkina Faso"
$ws.Range("B90").Value = 207
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 21
$ws.Range("E90").Value = 175
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 11

$ws.Range("A91").Value = "Vietnam"
$ws.Range("B91").Value = 188
$ws.Range("C91").Value = 14
$ws.Range("D91").Value = 21
$ws.Range("E91").Value = 167
$ws.Range("F91").Value = 3
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0

$ws.Range("A92").Value = "Reunion"
$ws.Range("B92").Value = 183
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 1
$ws.Range("E92").Value = 182
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0

$ws.Range("A104").Value = "Sri Lanka"
$ws.Range("B104").Value = 117
$ws.Range("C104").Value = 4
$ws.Range("D104").Value = 11
$ws.Range("E104").Value = 105
$ws.Range("F104").Value = 5
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 1

$ws.Range("A105").Value = "Honduras"
$ws.Range("B105").Value = 110
$ws.Range("C105").Value = 15
$ws.Range("D105").Value = 3
$ws.Range("E105").Value = 105
$ws.Range("F105").Value = 4
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 2

$ws.Range("A106").Value = "Mauricio"
$ws.Range("B106").Value = 107
$ws.Range("C106").Value = 5
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 105
$ws.Range("F106").Value = 1
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 2

$ws.Range("A107").Value = "Estado de Palestina"
$ws.Range("B107").Value = 106
$ws.Range("C107").Value = 2
$ws.Range("D107").Value = 18
$ws.Range("E107").Value = 87
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 1

$ws.Range("A108").Value = "Camboya"
$ws.Range("B108").Value = 103
$ws.Range("C108").Value = 4
$ws.Range("D108").Value = 21
$ws.Range("E108").Value = 82
$ws.Range("F108").Value = 1
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 0

$ws.Range("A169").Value = "Guyana"
$ws.Range("B169").Value = 8
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 0
$ws.Range("E169").Value = 7
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 1

$ws.Range("A170").Value = "Islas Caimanes"
$ws.Range("B170").Value = 8
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 0
$ws.Range("E170").Value = 7
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 1

$ws.Range("A175").Value = "San Martin (Parte Holandesa)"
$ws.Range("B175").Value = 6
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 0
$ws.Range("E175").Value = 6
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

$ws.Range("A176").Value = "Eritrea"
$ws.Range("B176").Value = 6
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 6
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

$ws.Range("A177").Value = "Santa Sede"
$ws.Range("B177").Value = 6
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 0
$ws.Range("E177").Value = 6
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

$ws.Range("A178").Value = "Benin"
$ws.Range("B178").Value = 6
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 6
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

$ws.Range("A179").Value = "Cabo Verde"
$ws.Range("B179").Value = 6
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 5
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 1

$ws.Range("A180").Value = "Fiyi"
$ws.Range("B180").Value = 5
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 5
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 0

$ws.Range("A181").Value = "Siria"
$ws.Range("B181").Value = 5
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 5
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0

$ws.Range("A182").Value = "Angola"
$ws.Range("B182").Value = 5
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 5
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 0

$ws.Range("A183").Value = "Mauritania"
$ws.Range("B183").Value = 5
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 5
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

$ws.Range("A184").Value = "Montserrat"
$ws.Range("B184").Value = 5
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 5
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

$ws.Range("A185").Value = "San Bartolome"
$ws.Range("B185").Value = 5
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 5
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0

$ws.Range("A187").Value = "Sudan"
$ws.Range("B187").Value = 5
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 4
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 1

$ws.Range("A188").Value = "Islas Turcas y Caicos"
$ws.Range("B188").Value = 4
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 4
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

$ws.Range("A189").Value = "Butan"
$ws.Range("B189").Value = 4
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 4
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 0

$ws.Range("A191").Value = "Santa Lucia"
$ws.Range("B191").Value = 4
$ws.Range("C191").Value = 1
$ws.Range("D191").Value = 1
$ws.Range("E191").Value = 3
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 0

$ws.Range("A192").Value = "Somalia"
$ws.Range("B192").Value = 3
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 0
$ws.Range("E192").Value = 3
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0
